$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: add F23 and G23 with value 0 (previously missing cells)
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0

# Update existing AgTests (F) / AgPosit (G) values for rows 269-386
$ws.Cells.Item(269, 6).Value = 7692
$ws.Cells.Item(269, 7).Value = 362
$ws.Cells.Item(270, 6).Value = 2423
$ws.Cells.Item(270, 7).Value = 155
$ws.Cells.Item(271, 6).Value = 43499
$ws.Cells.Item(271, 7).Value = 1679
$ws.Cells.Item(272, 6).Value = 30595
$ws.Cells.Item(272, 7).Value = 1648
$ws.Cells.Item(273, 6).Value = 31428
$ws.Cells.Item(273, 7).Value = 1653
$ws.Cells.Item(274, 6).Value = 28023
$ws.Cells.Item(274, 7).Value = 1285
$ws.Cells.Item(275, 6).Value = 29741
$ws.Cells.Item(275, 7).Value = 1274
$ws.Cells.Item(276, 6).Value = 11569
$ws.Cells.Item(276, 7).Value = 392
$ws.Cells.Item(277, 6).Value = 3387
$ws.Cells.Item(277, 7).Value = 127
$ws.Cells.Item(278, 6).Value = 29857
$ws.Cells.Item(278, 7).Value = 2058
$ws.Cells.Item(279, 6).Value = 42831
$ws.Cells.Item(279, 7).Value = 3050
$ws.Cells.Item(280, 6).Value = 34434
$ws.Cells.Item(280, 7).Value = 2293
$ws.Cells.Item(281, 6).Value = 46231
$ws.Cells.Item(281, 7).Value = 3180
$ws.Cells.Item(282, 6).Value = 46403
$ws.Cells.Item(282, 7).Value = 2764
$ws.Cells.Item(283, 6).Value = 17181
$ws.Cells.Item(284, 6).Value = 1229
$ws.Cells.Item(284, 7).Value = 100
$ws.Cells.Item(285, 6).Value = 41950
$ws.Cells.Item(285, 7).Value = 3410
$ws.Cells.Item(286, 6).Value = 55144
$ws.Cells.Item(286, 7).Value = 4281
$ws.Cells.Item(287, 6).Value = 58805
$ws.Cells.Item(287, 7).Value = 3719
$ws.Cells.Item(288, 6).Value = 58883
$ws.Cells.Item(288, 7).Value = 3938
$ws.Cells.Item(289, 6).Value = 63659
$ws.Cells.Item(289, 7).Value = 3682
$ws.Cells.Item(290, 6).Value = 17585
$ws.Cells.Item(290, 7).Value = 1038
$ws.Cells.Item(291, 6).Value = 15113
$ws.Cells.Item(292, 6).Value = 82459
$ws.Cells.Item(292, 7).Value = 7275
$ws.Cells.Item(293, 6).Value = 82707
$ws.Cells.Item(293, 7).Value = 5763
$ws.Cells.Item(294, 6).Value = 93882
$ws.Cells.Item(294, 7).Value = 4948
$ws.Cells.Item(295, 6).Value = 17112
$ws.Cells.Item(295, 7).Value = 1032
$ws.Cells.Item(296, 6).Value = 1874
$ws.Cells.Item(296, 7).Value = 142
$ws.Cells.Item(297, 6).Value = 2316
$ws.Cells.Item(298, 6).Value = 3211
$ws.Cells.Item(299, 6).Value = 65675
$ws.Cells.Item(299, 7).Value = 6877
$ws.Cells.Item(300, 6).Value = 72497
$ws.Cells.Item(300, 7).Value = 6976
$ws.Cells.Item(301, 6).Value = 72103
$ws.Cells.Item(301, 7).Value = 5681
$ws.Cells.Item(302, 6).Value = 78581
$ws.Cells.Item(302, 7).Value = 5655
$ws.Cells.Item(303, 6).Value = 9597
$ws.Cells.Item(303, 7).Value = 613
$ws.Cells.Item(304, 6).Value = 6052
$ws.Cells.Item(304, 7).Value = 521
$ws.Cells.Item(305, 6).Value = 3378
$ws.Cells.Item(305, 7).Value = 263
$ws.Cells.Item(306, 6).Value = 73994
$ws.Cells.Item(306, 7).Value = 7476
$ws.Cells.Item(307, 6).Value = 74892
$ws.Cells.Item(307, 7).Value = 6283
$ws.Cells.Item(308, 6).Value = 15471
$ws.Cells.Item(308, 7).Value = 1050
$ws.Cells.Item(309, 6).Value = 77458
$ws.Cells.Item(309, 7).Value = 5455
$ws.Cells.Item(310, 6).Value = 79199
$ws.Cells.Item(310, 7).Value = 4066
$ws.Cells.Item(311, 6).Value = 61506
$ws.Cells.Item(311, 7).Value = 1928
$ws.Cells.Item(312, 6).Value = 28135
$ws.Cells.Item(312, 7).Value = 926
$ws.Cells.Item(313, 6).Value = 75401
$ws.Cells.Item(313, 7).Value = 3448
$ws.Cells.Item(314, 6).Value = 64284
$ws.Cells.Item(315, 6).Value = 56370
$ws.Cells.Item(315, 7).Value = 2628
$ws.Cells.Item(316, 6).Value = 50720
$ws.Cells.Item(317, 6).Value = 63704
$ws.Cells.Item(317, 7).Value = 2173
$ws.Cells.Item(318, 6).Value = 48964
$ws.Cells.Item(318, 7).Value = 1135
$ws.Cells.Item(319, 6).Value = 41322
$ws.Cells.Item(320, 6).Value = 71356
$ws.Cells.Item(320, 7).Value = 3302
$ws.Cells.Item(321, 6).Value = 89210
$ws.Cells.Item(322, 6).Value = 109841
$ws.Cells.Item(322, 7).Value = 2348
$ws.Cells.Item(323, 6).Value = 216787
$ws.Cells.Item(324, 6).Value = 239328
$ws.Cells.Item(324, 7).Value = 2762
$ws.Cells.Item(325, 6).Value = 766053
$ws.Cells.Item(325, 7).Value = 6463
$ws.Cells.Item(327, 6).Value = 224563
$ws.Cells.Item(327, 7).Value = 2712
$ws.Cells.Item(328, 6).Value = 180496
$ws.Cells.Item(328, 7).Value = 2646
$ws.Cells.Item(329, 6).Value = 82715
$ws.Cells.Item(329, 7).Value = 1726
$ws.Cells.Item(330, 6).Value = 72515
$ws.Cells.Item(330, 7).Value = 2083
$ws.Cells.Item(331, 6).Value = 154628
$ws.Cells.Item(332, 6).Value = 457472
$ws.Cells.Item(332, 7).Value = 4540
$ws.Cells.Item(333, 6).Value = 271708
$ws.Cells.Item(334, 6).Value = 196216
$ws.Cells.Item(334, 7).Value = 3459
$ws.Cells.Item(335, 6).Value = 130351
$ws.Cells.Item(335, 7).Value = 2963
$ws.Cells.Item(336, 6).Value = 101484
$ws.Cells.Item(336, 7).Value = 3324
$ws.Cells.Item(337, 6).Value = 103581
$ws.Cells.Item(337, 7).Value = 2896
$ws.Cells.Item(338, 6).Value = 226616
$ws.Cells.Item(338, 7).Value = 3158
$ws.Cells.Item(339, 6).Value = 660435
$ws.Cells.Item(339, 7).Value = 5496
$ws.Cells.Item(341, 6).Value = 290886
$ws.Cells.Item(341, 7).Value = 3637
$ws.Cells.Item(342, 6).Value = 177063
$ws.Cells.Item(342, 7).Value = 2986
$ws.Cells.Item(343, 6).Value = 132607
$ws.Cells.Item(344, 6).Value = 135585
$ws.Cells.Item(344, 7).Value = 2486
$ws.Cells.Item(345, 6).Value = 291241
$ws.Cells.Item(345, 7).Value = 3306
$ws.Cells.Item(346, 6).Value = 675239
$ws.Cells.Item(346, 7).Value = 4837
$ws.Cells.Item(347, 6).Value = 344352
$ws.Cells.Item(347, 7).Value = 2927
$ws.Cells.Item(348, 6).Value = 231443
$ws.Cells.Item(348, 7).Value = 3222
$ws.Cells.Item(349, 6).Value = 158923
$ws.Cells.Item(349, 7).Value = 2752
$ws.Cells.Item(350, 6).Value = 127266
$ws.Cells.Item(350, 7).Value = 2782
$ws.Cells.Item(351, 6).Value = 150610
$ws.Cells.Item(352, 6).Value = 306504
$ws.Cells.Item(352, 7).Value = 3534
$ws.Cells.Item(353, 6).Value = 724705
$ws.Cells.Item(353, 7).Value = 5291
$ws.Cells.Item(354, 6).Value = 312433
$ws.Cells.Item(354, 7).Value = 2859
$ws.Cells.Item(355, 6).Value = 221687
$ws.Cells.Item(355, 7).Value = 3446
$ws.Cells.Item(356, 6).Value = 159809
$ws.Cells.Item(356, 7).Value = 2875
$ws.Cells.Item(357, 6).Value = 138228
$ws.Cells.Item(357, 7).Value = 3021
$ws.Cells.Item(358, 6).Value = 157305
$ws.Cells.Item(359, 6).Value = 320436
$ws.Cells.Item(359, 7).Value = 3346
$ws.Cells.Item(360, 6).Value = 747907
$ws.Cells.Item(360, 7).Value = 5126
$ws.Cells.Item(362, 6).Value = 227849
$ws.Cells.Item(362, 7).Value = 3167
$ws.Cells.Item(363, 6).Value = 188267
$ws.Cells.Item(363, 7).Value = 2766
$ws.Cells.Item(364, 6).Value = 167049
$ws.Cells.Item(364, 7).Value = 2455
$ws.Cells.Item(365, 6).Value = 183147
$ws.Cells.Item(365, 7).Value = 2389
$ws.Cells.Item(366, 6).Value = 338721
$ws.Cells.Item(367, 6).Value = 762617
$ws.Cells.Item(367, 7).Value = 3896
$ws.Cells.Item(368, 6).Value = 345571
$ws.Cells.Item(369, 6).Value = 233543
$ws.Cells.Item(369, 7).Value = 2590
$ws.Cells.Item(370, 6).Value = 181260
$ws.Cells.Item(370, 7).Value = 2019
$ws.Cells.Item(371, 6).Value = 158171
$ws.Cells.Item(371, 7).Value = 1941
$ws.Cells.Item(372, 6).Value = 177283
$ws.Cells.Item(372, 7).Value = 1845
$ws.Cells.Item(373, 6).Value = 346557
$ws.Cells.Item(374, 6).Value = 768046
$ws.Cells.Item(374, 7).Value = 3401
$ws.Cells.Item(375, 6).Value = 351563
$ws.Cells.Item(376, 6).Value = 220187
$ws.Cells.Item(376, 7).Value = 2204
$ws.Cells.Item(377, 6).Value = 175494
$ws.Cells.Item(377, 7).Value = 1798
$ws.Cells.Item(378, 6).Value = 155819
$ws.Cells.Item(378, 7).Value = 1524
$ws.Cells.Item(379, 6).Value = 178093
$ws.Cells.Item(379, 7).Value = 1600
$ws.Cells.Item(380, 6).Value = 341213
$ws.Cells.Item(380, 7).Value = 1983
$ws.Cells.Item(381, 6).Value = 737663
$ws.Cells.Item(381, 7).Value = 2648
$ws.Cells.Item(382, 6).Value = 355862
$ws.Cells.Item(382, 7).Value = 1561
$ws.Cells.Item(383, 6).Value = 218273
$ws.Cells.Item(383, 7).Value = 1738
$ws.Cells.Item(384, 6).Value = 167699
$ws.Cells.Item(384, 7).Value = 1475
$ws.Cells.Item(385, 6).Value = 145497
$ws.Cells.Item(385, 7).Value = 1365
$ws.Cells.Item(386, 6).Value = 176037
$ws.Cells.Item(386, 7).Value = 1323

# Append new rows 387 and 388
$ws.Cells.Item(387, 1).Value = 44281
$ws.Cells.Item(387, 2).Value = 356985
$ws.Cells.Item(387, 3).Value = 12389
$ws.Cells.Item(387, 4).Value = 1531
$ws.Cells.Item(387, 5).Value = 9426
$ws.Cells.Item(387, 6).Value = 325997
$ws.Cells.Item(387, 7).Value = 1530

$ws.Cells.Item(388, 1).Value = 44282
$ws.Cells.Item(388, 2).Value = 357910
$ws.Cells.Item(388, 3).Value = 7087
$ws.Cells.Item(388, 4).Value = 925
$ws.Cells.Item(388, 5).Value = 9496
$ws.Cells.Item(388, 6).Value = 592025
$ws.Cells.Item(388, 7).Value = 2397
